$d = $word.ActiveDocument

# The page ends with:
#   ... "LOB1012: Estatística (Requisito)"
#   <empty paragraph>
#   <empty paragraph with pageBreakBefore>
#   "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#    pages. Original theme under Creative Commons Attribution"
#   <empty paragraph>
#   <empty paragraph with pageBreakBefore>
#
# The footer block (the blank/page-break paragraphs plus the copyright
# paragraph that immediately follow the "LOB1012..." requirement line)
# was removed, leaving only the trailing blank + page-break paragraphs.
# Locate the anchor paragraphs by their text instead of hard-coded
# indices so the edit is resilient to any earlier content differences.

$anchorIndex = -1
$copyrightIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($t -like "*LOB1012*") {
        $anchorIndex = $i
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $copyrightIndex = $i
    }
}

if ($anchorIndex -ge 0 -and $copyrightIndex -ge $anchorIndex) {
    $startPara = $d.Paragraphs.Item($anchorIndex + 1)
    $endPara = $d.Paragraphs.Item($copyrightIndex)

    $start = $startPara.Range.Start
    $end = $endPara.Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
